$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I (I0) and J (IF), reusing the
# existing header style (same as H1) via a format-only paste so the
# shared style index is reused instead of creating a near-duplicate.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows 2-72 for columns I (I0) and J (IF), as parallel arrays
# (row number / I value / J value at the same index).
$rowNums = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72)
$iVals   = @(9, 8, 7, 8, 9, 9, 8, 9, 9, 9, 7, 8, 9, 8, 7, 9, 8, 9, 8, 8, 8, 8, 8, 6, 8, 8, 8, 7, 9, 8, 7, 8, 7, 7, 6, 7, 9, 8, 8, 7, 7, 7, 7, 9, 9, 8, 8, 8, 7, 8, 7, 9, 5, 6, 9, 6, 9, 8, 8, 8, 7, 3, 7, 4, 9, 7, 4, 7, 9, 5, 4)
$jVals   = @(9, 8, 7, 8, 9, 10, 8, 9, 9, 9, 8, 8, 9, 9, 8, 10, 8, 9, 8, 8, 8, 8, 8, 6, 8, 8, 8, 7, 9, 9, 7, 8, 7, 7, 6, 7, 9, 8, 8, 7, 7, 7, 7, 9, 9, 8, 8, 8, 7, 8, 8, 9, 5, 6, 9, 7, 9, 8, 8, 8, 7, 3, 7, 4, 9, 8, 4, 7, 9, 5, 4)

for ($idx = 0; $idx -lt $rowNums.Length; $idx++) {
    $r = $rowNums[$idx]
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}

